# Update gh-pages to output generated at 456a3b4
# Bumps a handful of view-count numbers (column F) across the
# "展览" (Exhibitions), "演出" (Shows) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 56
$ws1.Range("F8").Value  = 53
$ws1.Range("F14").Value = 2016
$ws1.Range("F23").Value = 1513
$ws1.Range("F24").Value = 3426
$ws1.Range("F35").Value = 405
$ws1.Range("F39").Value = 372

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 56
$ws4.Range("F8").Value  = 53
$ws4.Range("F14").Value = 2016
$ws4.Range("F16").Value = 12
$ws4.Range("F24").Value = 1513
$ws4.Range("F25").Value = 3426
$ws4.Range("F36").Value = 405
$ws4.Range("F40").Value = 373
